# "completed part 3 and 4"
#
# The workbook previously reserved column A for an "Audio" file-name
# reference (header "Audio" / value "audio.mp3") with the question's
# answer-key rows (A2:A4) merged together. That column is no longer
# needed, so it is removed entirely - every other column (B:G) shifts
# one place to the left (becoming A:F) and the merged placeholder cell
# disappears along with the deleted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A (the old "Audio"/"audio.mp3" column). Excel shifts the
# remaining columns left, drops the now-fully-deleted A2:A4 merge, and
# shrinks the sheet's dimension/used range accordingly.
$ws.Columns.Item(1).Delete()

# Leave the selection on column A (now holding what used to be column B),
# matching the saved view state.
$null = $ws.Columns.Item(1).Select()
